# Generate Report for Handback
#
# The f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md file has now been handed
# back (for both the zh-cn and de-de targets), so update the localization
# status report: the "Ready for handoff" status becomes
# "Handed back: in sync with en-US", the Latest Handback DateTime is
# refreshed, and the stale "version ... is not the latest" error is
# cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
# Row 3 is the f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md file; its zh-cn
# (E) and de-de (F) status columns move from "Ready for handoff" to
# "Handed back: in sync with en-US".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
# Row 3 is the f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md handback row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-13 21:03:37"
$wsZhCn.Range("P3").Value = ""
# Error Detail column shrinks now that it no longer holds the long stale-handback message.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet --------------------------------------------------------
# Row 3 is the f5fcffff-7c7b-4436-a1ba-0c7e540837c7.md handback row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-13 21:03:47"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
